# Updates cryptos list values (Price / Volume(1h) columns) and a couple of
# row reorderings (Coin name + Link) to match the refreshed data snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay plain text even when it looks like a
# number (e.g. "0.430", "1.85") so Excel does not silently coerce it to a
# numeric cell (which would drop trailing zeros / change the cell type).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "60.996.02"
$ws.Range("E2").Value = "  +1.25%  "
$ws.Range("D3").Value = "2.632.07"
$ws.Range("E3").Value = "  +1.66%  "
$ws.Range("E4").Value = "  -0.05%  "
Set-TextValue $ws.Range("D5") "529.95"
$ws.Range("E5").Value = "  +4.17%  "
Set-TextValue $ws.Range("D6") "155.53"
$ws.Range("E6").Value = "  +1.17%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.09%  "
Set-TextValue $ws.Range("D9") "6.67"
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("E10").Value = "  +5.66%  "
$ws.Range("E11").Value = "  +1.06%  "
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("D13").Value = "3.097.96"
$ws.Range("E13").Value = "  +1.71%  "
$ws.Range("D14").Value = "61.000.00"
$ws.Range("E14").Value = "  +1.24%  "
Set-TextValue $ws.Range("D15") "22.04"
$ws.Range("E15").Value = "  +2.39%  "
$ws.Range("E16").Value = "  +3.65%  "
$ws.Range("D17").Value = "2.638.01"
$ws.Range("E17").Value = "  +1.64%  "
Set-TextValue $ws.Range("D18") "4.77"
$ws.Range("E18").Value = "  +0.68%  "
Set-TextValue $ws.Range("D19") "356.60"
$ws.Range("E19").Value = "  +1.34%  "
Set-TextValue $ws.Range("D20") "10.64"
$ws.Range("E20").Value = "  +0.99%  "
Set-TextValue $ws.Range("D21") "6.25"
$ws.Range("E21").Value = "  +2.24%  "
$ws.Range("E22").Value = "  -0.02%  "
Set-TextValue $ws.Range("D23") "61.81"
$ws.Range("E23").Value = "  +2.36%  "
Set-TextValue $ws.Range("D24") "0.430"
$ws.Range("E24").Value = "  +2.35%  "
Set-TextValue $ws.Range("D25") "0.169"
$ws.Range("E25").Value = "  +1.44%  "
Set-TextValue $ws.Range("D26") "0.999"
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("D27").Value = "0.0₃0871"
$ws.Range("E27").Value = "  +3.90%  "
Set-TextValue $ws.Range("D28") "7.42"
$ws.Range("E28").Value = "  +1.17%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D30") "6.15"
$ws.Range("E30").Value = "  +7.34%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D31") "19.52"
$ws.Range("E31").Value = "  +0.70%  "
$ws.Range("E32").Value = "  +4.23%  "
Set-TextValue $ws.Range("D33") "151.37"
$ws.Range("E33").Value = "  -0.37%  "
Set-TextValue $ws.Range("D34") "4.16"
$ws.Range("E34").Value = "  +4.23%  "
$ws.Range("E35").Value = "  +1.92%  "
Set-TextValue $ws.Range("D36") "0.934"
$ws.Range("E36").Value = "  +11.33%  "
Set-TextValue $ws.Range("D37") "0.885"
$ws.Range("E37").Value = "  +3.03%  "
$ws.Range("E38").Value = "  +1.52%  "
Set-TextValue $ws.Range("D39") "3.82"
$ws.Range("E39").Value = "  +1.89%  "
Set-TextValue $ws.Range("D40") "301.37"
$ws.Range("E40").Value = "  +1.92%  "
Set-TextValue $ws.Range("D41") "0.642"
$ws.Range("E41").Value = "  +3.80%  "
$ws.Range("E42").Value = "  +1.50%  "
$ws.Range("E43").Value = "  +1.85%  "
Set-TextValue $ws.Range("D44") "0.998"
$ws.Range("E44").Value = "  -0.01%  "
Set-TextValue $ws.Range("D45") "19.75"
$ws.Range("E45").Value = "  +0.54%  "
Set-TextValue $ws.Range("D46") "4.99"
$ws.Range("E46").Value = "  +4.21%  "
$ws.Range("E47").Value = "  +2.85%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D48") "19.24"
$ws.Range("E48").Value = "  +8.14%  "
$ws.Range("B49").Value = "WhiteBITCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue $ws.Range("D49") "10.35"
$ws.Range("E49").Value = "  +0.49%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "1.986.65"
$ws.Range("E50").Value = "  -0.22%  "
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D51") "1.85"
$ws.Range("E51").Value = "  +3.85%  "
